$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -4
    3  = -1
    4  = 2
    5  = -1
    6  = -4
    7  = -2
    8  = -3
    9  = -8
    10 = -3
    11 = -3
    12 = -1
    13 = 0
    14 = -5
    15 = 4
    16 = -3
    17 = -2
    18 = -5
    20 = 1
    22 = -5
    23 = -3
    24 = 0
    25 = 2
    26 = 5
    28 = 2
    29 = -3
    30 = 5
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
